# "Fruta / hortaliza, semanal" -- weekly refresh of the Jengibre (ginger)
# subconjunto: the price-observation rows (Fecha, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Precio $/Kg) get reshuffled
# across the data rows (2-45) of the sheet. Row 17 is untouched. Every
# other row receives the D/J/K/L/M/P tuple that used to live on a
# (different) row, per the mapping below -- all other columns (A, B, C,
# E, F, G, H, I, N, O, Q, R) are identical on every row already, so they
# don't need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row -> source row (the row whose D/J/K/L/M/P values now land here)
$mapping = @{
    2  = 33
    3  = 27
    4  = 18
    5  = 30
    6  = 12
    7  = 20
    8  = 19
    9  = 23
    10 = 36
    11 = 21
    12 = 15
    13 = 40
    14 = 2
    15 = 32
    16 = 9
    17 = 17
    18 = 11
    19 = 16
    20 = 8
    21 = 45
    22 = 26
    23 = 44
    24 = 28
    25 = 34
    26 = 13
    27 = 24
    28 = 5
    29 = 7
    30 = 10
    31 = 4
    32 = 39
    33 = 31
    34 = 14
    35 = 29
    36 = 43
    37 = 35
    38 = 22
    39 = 37
    40 = 38
    41 = 6
    42 = 3
    43 = 25
    44 = 41
    45 = 42
}

$cols = @("D", "J", "K", "L", "M", "P")

# Snapshot every source row's current values before writing anything,
# so overlapping source/target rows never read back an already-updated
# value.
$snapshot = @{}
foreach ($row in 2..45) {
    $values = @{}
    foreach ($col in $cols) {
        $values[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $values
}

foreach ($row in 2..45) {
    $srcRow = $mapping[$row]
    if ($srcRow -eq $row) {
        continue
    }
    $values = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
